$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.945.58'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.844.48'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '692.60'
$ws.Range("E5").Value = '  +3.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.80'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("D7").Value = '3.842.76'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.38'
$ws.Range("E11").Value = '  +5.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("E13").Value = '  +5.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.61'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("D15").Value = '4.510.34'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '3.870.70'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '71.037.95'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.82'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.24'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.18'
$ws.Range("E21").Value = '  -4.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.81'
$ws.Range("E22").Value = '  +2.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.70'
$ws.Range("E24").Value = '  +1.76%  '
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.39'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.49'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.14'
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = '4.003.70'
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +8.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.61'
$ws.Range("E32").Value = '  +2.51%  '
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.70'
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.25'
$ws.Range("E36").Value = '  +0.76%  '
$ws.Range("D37").Value = '3.797.46'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.38'
$ws.Range("E40").Value = '  +13.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.44'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.04'
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("E43").Value = '  +4.19%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '164.51'
$ws.Range("E46").Value = '  +4.45%  '
$ws.Range("E47").Value = '  +5.80%  '
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.55'
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.39'
$ws.Range("E51").Value = '  -2.74%  '
